$wb = $excel.ActiveWorkbook

# --- TestSteps sheet: update the TestData-column references to use the new prefixed names ---
$ws2 = $wb.Worksheets.Item("TestSteps")
$ws2.Range("F2").Value = "d_browser"
$ws2.Range("F4").Value = "d_username"
$ws2.Range("F5").Value = "d_password"

# --- TestData sheet: prefix the testData column headers with "d_" ---
$ws3 = $wb.Worksheets.Item("TestData")
$ws3.Range("B1").Value = "d_username"
$ws3.Range("C1").Value = "d_password"
$ws3.Range("D1").Value = "d_browser"
$ws3.Range("B1:D1").EntireColumn.AutoFit()

# --- Update selection on TestSteps (no longer the active tab) ---
[void]$ws2.Activate()
[void]$ws2.Range("F5").Select()

# --- Make TestData the active tab with its new selection ---
[void]$ws3.Activate()
[void]$ws3.Range("D9").Select()
